$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text block in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("✅ 1000 Bs = 12.36 = 49753.03 pesos", "✅ 1000 Bs = 12.64 = 50941.85 pesos")
$text = $text.Replace("✅ 49753.03 pesos = 12.35 = 967.98 Bs", "✅ 50941.85 pesos = 12.61 = 958.31 Bs")
$cell.Value = $text

# --- Sheet "tasas": update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 79.11
$wsTasas.Range("O10").Value = 4030.01
$wsTasas.Range("N12").Value = 4040
$wsTasas.Range("O12").Value = 76
